$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Value = "Riccardo Barbiero"
$ws.Range("B67").Value = "Elia Battisti | U.SGUARNA"
$ws.Range("C67").Value = "Riccardo Barbiero | Rita Levi’s"
$ws.Range("D67").Value = "Leonardo Viola | SHARK ATTACK"
$ws.Range("E67").Value = "Marco Sala | IMONTAGNA"
$ws.Range("F67").Value = "Moris Benedetti | Gli Introvabili"
